# Commit: "macro formulated cells are modified for debugging in date
# conversion"
#
# The Aircraft_scheduling sheet computed the work-package "finish time"
# (column E) from the "start time" (column C) with a formula such as
# `=C6+0.181`. While debugging a date/time-rollover issue, those formula
# cells were replaced with their plain calculated numbers (the shared
# MOD(E-C,1) duration formula in column F is left in place and simply
# recalculates against the new literal values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

$ws.Range("E6").Value  = 0.6462731481481482
$ws.Range("E11").Value = 0.2469675925925926
$ws.Range("E12").Value = 0.1322222222222222
$ws.Range("E13").Value = 0.33375
$ws.Range("E14").Value = 0.4306365740740741
$ws.Range("E15").Value = 0.26055555555555554
$ws.Range("E16").Value = 0.29185185185185186
$ws.Range("E17").Value = 0.3616898148148148
$ws.Range("E18").Value = 0.5489351851851852

# Match the author's final on-screen selection: the Aircraft_scheduling
# tab became the active/selected sheet (moving tabSelected away from the
# Staff sheet) with E16 highlighted.
$ws.Select()
$ws.Range("E16").Select()
